# Updated model dataframes to exclude outliers.
# This rewrites the data rows (category name, positives, ranks, pos_ranked)
# on both the "covariate_importance" and "strategy_importance" sheets with
# refreshed statistics (outliers excluded) and re-sorts them descending by
# the "positives" column, matching the new 03identify_outliers.R pipeline.

$wb = $excel.ActiveWorkbook

$covariateData = @(
    @("region", 100, 100, 100),
    @("state", 100, 100, 100),
    @("percenttwoormoreraces", 87, 57.99999999999999, 57.99999999999999),
    @("percentwhite", 84, 50, 50),
    @("rplthemes", 69, 42, 42),
    @("percentblackorafricanamerican", 75, 40, 40),
    @("percentasian", 80, 39, 39),
    @("percentstudentsfreereducedlunch", 74, 36, 36),
    @("percentamericanindianoralaskanative", 45, 9, 9),
    @("cntycaseschange", 41, 6, 6),
    @("percenthispaniclatino", 34, 6, 6),
    @("derivedtotalenrolled", 39, 5, 5),
    @("percentnativehawaiianorotherpacificislander", 23, 5, 5),
    @("schoollevel", 37, 4, 4),
    @("locale", 14, 0, 0),
    @("percentnotspecified", 0, 0, 0)
)

$strategyData = @(
    @("cleaning", 100, 100, 100),
    @("physicaldistancing", 65, 73, 63),
    @("quarantine", 62, 72, 62),
    @("contacttracing", 32, 38, 31),
    @("hvacsystems", 29, 37, 28),
    @("vaccination", 28, 31, 26),
    @("masks", 22, 24, 22),
    @("screeningtestingforstudents", 8, 11, 7.000000000000001),
    @("hepafilters", 0, 14, 0),
    @("stayhome", 0, 0, 0)
)

$wsCovariate = $wb.Worksheets.Item("covariate_importance")
$row = 2
foreach ($entry in $covariateData) {
    $wsCovariate.Cells.Item($row, 1).Value = $entry[0]
    $wsCovariate.Cells.Item($row, 2).Value = $entry[1]
    $wsCovariate.Cells.Item($row, 3).Value = $entry[2]
    $wsCovariate.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}

$wsStrategy = $wb.Worksheets.Item("strategy_importance")
$row = 2
foreach ($entry in $strategyData) {
    $wsStrategy.Cells.Item($row, 1).Value = $entry[0]
    $wsStrategy.Cells.Item($row, 2).Value = $entry[1]
    $wsStrategy.Cells.Item($row, 3).Value = $entry[2]
    $wsStrategy.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
